$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows with new Esperado/Observado/valor p values (semana 41 de 2025)
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = 0.06

$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 0

$ws.Range("C5").Value = 6
$ws.Range("D5").Value = 13
$ws.Range("E5").Value = 0.01

$ws.Range("C6").Value = 3
$ws.Range("E6").Value = 0.22

$ws.Range("D7").Value = 2
$ws.Range("E7").Value = 0.27

$ws.Range("C9").Value = 48
$ws.Range("D9").Value = 41
$ws.Range("E9").Value = 0.04

$ws.Range("C10").Value = 1
$ws.Range("E10").Value = 0.37

$ws.Range("D11").Value = 3
$ws.Range("E11").Value = 0.18

$ws.Range("C12").Value = 5
$ws.Range("E12").Value = 0.18

$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 1

$ws.Range("C16").Value = 10
$ws.Range("E16").Value = 0

$ws.Range("C17").Value = 9
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = 0.13

$ws.Range("C18").Value = 2
$ws.Range("E18").Value = 0.14

$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 0.1

$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0.14

$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 0.08

$ws.Range("C26").Value = 0
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = 0

$ws.Range("C31").Value = 1
$ws.Range("D31").Value = 1
$ws.Range("E31").Value = 0.37

$ws.Range("C33").Value = 6
$ws.Range("D33").Value = 7
$ws.Range("E33").Value = 0.14

$ws.Range("D34").Value = 2
$ws.Range("E34").Value = 0

$ws.Range("C35").Value = 8
$ws.Range("D35").Value = 7
$ws.Range("E35").Value = 0.14

# Add new row 36 for event 895 - Zika
$ws.Range("A36").NumberFormat = "@"
$ws.Range("A36").Value = "895"
$ws.Range("B36").Value = "Zika"
$ws.Range("C36").Value = 0
$ws.Range("D36").Value = 1
$ws.Range("E36").Value = 0
